$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: 05/31/2020 | 2 | Figures: productive linkages task 3.1 (complete) and 3.2 (in progress)
# A9 must stay literal text (not get auto-parsed into a date serial), matching the
# look of the existing "date-like" text cells above it (A5-A8).
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "05/31/2020"
$ws.Range("A9").ClearFormats()

$ws.Range("B9").Value = 2

$ws.Range("C9").Value = "Figures: productive linkages task 3.1 (complete) and 3.2 (in progress)"
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Rows.Item(9).RowHeight = 15.75

# --- Row 10: 05/31/2020 (serial 43836) | 3 | Figures: productive linkages part 3.1 and 3.2 initial proposal
$ws.Range("A10").Value = 43836
$ws.Range("A10").NumberFormat = "m/d/yy"

$ws.Range("B10").Value = 3

$ws.Range("C10").Value = "Figures: productive linkages part 3.1 and 3.2 initial proposal"
$ws.Range("C8").Copy()
$ws.Range("C10").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Rows.Item(10).RowHeight = 15.75

$excel.CutCopyMode = $false

$ws.Range("A10").Select()
